# Scheduled-runner price refresh: updates computed Leve profit columns (H-N)
# on several rows across multiple class sheets. Values below are taken from
# the freshly recalculated market data for this run.
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# --- ALC ---
# row 18
$ws_ALC.Range("H18").Value = 5185.2856
$ws_ALC.Range("I18").Value = 859.8
$ws_ALC.Range("K18").Value = 859.8
$ws_ALC.Range("M18").Value = -575.8
# row 42
$ws_ALC.Range("H42").Value = 242.25
$ws_ALC.Range("I42").Value = 169.78572
$ws_ALC.Range("J42").Value = 749.5
$ws_ALC.Range("K42").Value = 509.35716
$ws_ALC.Range("L42").Value = 2248.5
$ws_ALC.Range("M42").Value = -279.35716
$ws_ALC.Range("N42").Value = -2708.5
# row 87
$ws_ALC.Range("H87").Value = 93798
$ws_ALC.Range("J87").Value = 94747.5
$ws_ALC.Range("L87").Value = 94747.5
$ws_ALC.Range("N87").Value = -97243.5
# row 90
$ws_ALC.Range("H90").Value = 93798
$ws_ALC.Range("J90").Value = 94747.5
$ws_ALC.Range("L90").Value = 284242.5
$ws_ALC.Range("N90").Value = -296722.5
# row 98
$ws_ALC.Range("H98").Value = 2164.1667
$ws_ALC.Range("I98").Value = 2164
$ws_ALC.Range("J98").Value = 2164.6667
$ws_ALC.Range("K98").Value = 2164
$ws_ALC.Range("L98").Value = 2164.6667
$ws_ALC.Range("M98").Value = -666
$ws_ALC.Range("N98").Value = -5160.6667
# row 100
$ws_ALC.Range("H100").Value = 6211.706
$ws_ALC.Range("J100").Value = 9502.111000000001
$ws_ALC.Range("L100").Value = 9502.111000000001
$ws_ALC.Range("N100").Value = -10584.111
# row 122
$ws_ALC.Range("H122").Value = 2164.1667
$ws_ALC.Range("I122").Value = 2164
$ws_ALC.Range("J122").Value = 2164.6667
$ws_ALC.Range("K122").Value = 6492
$ws_ALC.Range("L122").Value = 6494.000100000001
$ws_ALC.Range("M122").Value = -4042
$ws_ALC.Range("N122").Value = -11394.0001
# row 137
$ws_ALC.Range("H137").Value = 34419900
$ws_ALC.Range("I137").Value = 47621164
$ws_ALC.Range("K137").Value = 142863492
$ws_ALC.Range("M137").Value = -142860942
# --- ARM ---
# row 28
$ws_ARM.Range("H28").Value = 1358.3334
$ws_ARM.Range("I28").Value = 1358.3334
$ws_ARM.Range("K28").Value = 1358.3334
$ws_ARM.Range("M28").Value = -1166.3334
# row 37
$ws_ARM.Range("H37").Value = 62423.105
$ws_ARM.Range("I37").Value = 22588.2
$ws_ARM.Range("J37").Value = 76649.86
$ws_ARM.Range("K37").Value = 22588.2
$ws_ARM.Range("L37").Value = 76649.86
$ws_ARM.Range("M37").Value = -22315.2
$ws_ARM.Range("N37").Value = -77195.86
# row 55
$ws_ARM.Range("H55").Value = 99977.5
$ws_ARM.Range("J55").Value = 99977.5
$ws_ARM.Range("L55").Value = 99977.5
$ws_ARM.Range("N55").Value = -100607.5
# row 61
$ws_ARM.Range("H61").Value = 4765907
$ws_ARM.Range("I61").Value = 5559391.5
$ws_ARM.Range("K61").Value = 5559391.5
$ws_ARM.Range("M61").Value = -5559179.5
# row 74
$ws_ARM.Range("H74").Value = 2376
$ws_ARM.Range("I74").Value = 859.05554
$ws_ARM.Range("K74").Value = 859.05554
$ws_ARM.Range("M74").Value = 14.94446000000005
# row 77
$ws_ARM.Range("H77").Value = 2376
$ws_ARM.Range("I77").Value = 859.05554
$ws_ARM.Range("K77").Value = 4295.2777
$ws_ARM.Range("M77").Value = 72.72230000000036
# row 80
$ws_ARM.Range("H80").Value = 89388.60000000001
$ws_ARM.Range("J80").Value = 92985.75
$ws_ARM.Range("L80").Value = 92985.75
$ws_ARM.Range("N80").Value = -94981.75
# row 83
$ws_ARM.Range("H83").Value = 89388.60000000001
$ws_ARM.Range("J83").Value = 92985.75
$ws_ARM.Range("L83").Value = 278957.25
$ws_ARM.Range("N83").Value = -288941.25
# row 88
$ws_ARM.Range("H88").Value = 3123.5
$ws_ARM.Range("J88").Value = 3214
$ws_ARM.Range("L88").Value = 3214
$ws_ARM.Range("N88").Value = -4026
# row 91
$ws_ARM.Range("H91").Value = 3123.5
$ws_ARM.Range("J91").Value = 3214
$ws_ARM.Range("L91").Value = 3214
$ws_ARM.Range("N91").Value = -6022
# row 99
$ws_ARM.Range("H99").Value = 1358.3334
$ws_ARM.Range("I99").Value = 1358.3334
$ws_ARM.Range("K99").Value = 1358.3334
$ws_ARM.Range("M99").Value = 1636.6666
# row 132
$ws_ARM.Range("H132").Value = 1669482.1
$ws_ARM.Range("I132").Value = 1820816.9
$ws_ARM.Range("K132").Value = 5462450.699999999
$ws_ARM.Range("M132").Value = -5459920.699999999
# row 136
$ws_ARM.Range("H136").Value = 4765907
$ws_ARM.Range("I136").Value = 5559391.5
$ws_ARM.Range("K136").Value = 16678174.5
$ws_ARM.Range("M136").Value = -16675624.5
# --- BSM ---
# row 86
$ws_BSM.Range("H86").Value = 1542.7142
$ws_BSM.Range("I86").Value = 1208.6
$ws_BSM.Range("K86").Value = 1208.6
$ws_BSM.Range("M86").Value = -85.59999999999991
# row 89
$ws_BSM.Range("H89").Value = 1542.7142
$ws_BSM.Range("I89").Value = 1208.6
$ws_BSM.Range("K89").Value = 6043
$ws_BSM.Range("M89").Value = -427
# row 94
$ws_BSM.Range("H94").Value = 1070.875
$ws_BSM.Range("I94").Value = 1081.7587
$ws_BSM.Range("K94").Value = 1081.7587
$ws_BSM.Range("M94").Value = -630.7587000000001
# row 134
$ws_BSM.Range("H134").Value = 1579053.5
$ws_BSM.Range("I134").Value = 1589473.8
$ws_BSM.Range("K134").Value = 4768421.4
$ws_BSM.Range("M134").Value = -4765886.4
# --- CRP ---
# row 17
$ws_CRP.Range("H17").Value = 3504.5
$ws_CRP.Range("I17").Value = 3504.5
$ws_CRP.Range("K17").Value = 3504.5
$ws_CRP.Range("M17").Value = -3330.5
# row 31
$ws_CRP.Range("H31").Value = 18954.424
$ws_CRP.Range("J31").Value = 54318.855
$ws_CRP.Range("L31").Value = 54318.855
$ws_CRP.Range("N31").Value = -54908.855
# row 34
$ws_CRP.Range("H34").Value = 18954.424
$ws_CRP.Range("J34").Value = 54318.855
$ws_CRP.Range("L34").Value = 54318.855
$ws_CRP.Range("N34").Value = -54722.855
# row 50
$ws_CRP.Range("H50").Value = 66499
$ws_CRP.Range("J50").Value = 112998
$ws_CRP.Range("L50").Value = 112998
$ws_CRP.Range("N50").Value = -114248
# row 51
$ws_CRP.Range("H51").Value = 25831.334
$ws_CRP.Range("J51").Value = 89988
$ws_CRP.Range("L51").Value = 89988
$ws_CRP.Range("N51").Value = -91460
# row 59
$ws_CRP.Range("H59").Value = 89997
$ws_CRP.Range("I59").Value = 30000
$ws_CRP.Range("K59").Value = 30000
$ws_CRP.Range("M59").Value = -28855
# row 60
$ws_CRP.Range("H60").Value = 30018
$ws_CRP.Range("J60").Value = 99998
$ws_CRP.Range("L60").Value = 99998
$ws_CRP.Range("N60").Value = -101020
# row 61
$ws_CRP.Range("H61").Value = 25831.334
$ws_CRP.Range("J61").Value = 89988
$ws_CRP.Range("L61").Value = 89988
$ws_CRP.Range("N61").Value = -90684
# row 68
$ws_CRP.Range("H68").Value = 100294.664
$ws_CRP.Range("J68").Value = 100294.664
$ws_CRP.Range("L68").Value = 100294.664
$ws_CRP.Range("N68").Value = -101792.664
# row 71
$ws_CRP.Range("H71").Value = 100294.664
$ws_CRP.Range("J71").Value = 100294.664
$ws_CRP.Range("L71").Value = 300883.992
$ws_CRP.Range("N71").Value = -308371.992
# row 74
$ws_CRP.Range("H74").Value = 64957
$ws_CRP.Range("J74").Value = 64957
$ws_CRP.Range("L74").Value = 64957
$ws_CRP.Range("N74").Value = -66705
# row 77
$ws_CRP.Range("H77").Value = 64957
$ws_CRP.Range("J77").Value = 64957
$ws_CRP.Range("L77").Value = 194871
$ws_CRP.Range("N77").Value = -203607
# row 86
$ws_CRP.Range("H86").Value = 339583.34
$ws_CRP.Range("I86").Value = 6885.75
$ws_CRP.Range("J86").Value = 1004978.5
$ws_CRP.Range("K86").Value = 6885.75
$ws_CRP.Range("L86").Value = 1004978.5
$ws_CRP.Range("M86").Value = -5762.75
$ws_CRP.Range("N86").Value = -1007224.5
# row 89
$ws_CRP.Range("H89").Value = 339583.34
$ws_CRP.Range("I89").Value = 6885.75
$ws_CRP.Range("J89").Value = 1004978.5
$ws_CRP.Range("K89").Value = 34428.75
$ws_CRP.Range("L89").Value = 5024892.5
$ws_CRP.Range("M89").Value = -28812.75
$ws_CRP.Range("N89").Value = -5036124.5
# --- CUL ---
# row 4
$ws_CUL.Range("H4").Value = 39626372
$ws_CUL.Range("I4").Value = 51389500
$ws_CUL.Range("J4").Value = 4337000
$ws_CUL.Range("K4").Value = 154168500
$ws_CUL.Range("L4").Value = 13011000
$ws_CUL.Range("M4").Value = -154168388
$ws_CUL.Range("N4").Value = -13011224
# row 12
$ws_CUL.Range("H12").Value = 30.307692
$ws_CUL.Range("J12").Value = 29.90909
$ws_CUL.Range("L12").Value = 89.72727
$ws_CUL.Range("N12").Value = -435.72727
# row 50
$ws_CUL.Range("H50").Value = 1799.8
$ws_CUL.Range("I50").Value = 333.33334
$ws_CUL.Range("K50").Value = 1000.00002
$ws_CUL.Range("M50").Value = -519.0000200000001
# row 53
$ws_CUL.Range("H53").Value = 1799.8
$ws_CUL.Range("I53").Value = 333.33334
$ws_CUL.Range("K53").Value = 1000.00002
$ws_CUL.Range("M53").Value = -519.0000200000001
# --- LTW ---
# row 103
$ws_LTW.Range("H103").Value = 82899.75
$ws_LTW.Range("J103").Value = 82899.75
$ws_LTW.Range("L103").Value = 82899.75
$ws_LTW.Range("N103").Value = -85243.75
# --- WVR ---
# row 126
$ws_WVR.Range("H126").Value = 2018.9474
$ws_WVR.Range("I126").Value = 2053.0557
$ws_WVR.Range("K126").Value = 6159.1671
$ws_WVR.Range("M126").Value = -3689.1671
# row 132
$ws_WVR.Range("H132").Value = 11804379
$ws_WVR.Range("J132").Value = 2199.5
$ws_WVR.Range("L132").Value = 6598.5
$ws_WVR.Range("N132").Value = -11658.5
# row 136
$ws_WVR.Range("H136").Value = 13899803
$ws_WVR.Range("I136").Value = 17555728
$ws_WVR.Range("K136").Value = 52667184
$ws_WVR.Range("M136").Value = -52664634
